$d = $word.ActiveDocument

# Helper characters
$rsquo = [char]0x2019   # curly apostrophe
$endash = [char]0x2013  # en dash

# ---------------------------------------------------------------
# 1) "Quelqu'un d'autre s'en charge"  (list item, ilvl = 1)
# ---------------------------------------------------------------
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.ListFormat.ListLevelNumber = 2
$p.Range.Text = "Quelqu" + $rsquo + "un d" + $rsquo + "autre s" + $rsquo + "en charge"

# ---------------------------------------------------------------
# 2) "Aide sur le core, probleme au demarrage" (list item, ilvl = 0)
# ---------------------------------------------------------------
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.ListFormat.ListLevelNumber = 1
$p.Range.Text = "Aide sur le core, probleme au demarrage"

# ---------------------------------------------------------------
# 3) empty paragraph
# ---------------------------------------------------------------
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.ListFormat.RemoveNumbers()
$p.Style = "Normal"

# ---------------------------------------------------------------
# 4) "Semaine 15" centered / bold / 16pt  +  5) date, centered
#    (insert the date placeholder first, then the heading
#    *before* it so the heading's bold formatting never leaks
#    forward onto the date line)
# ---------------------------------------------------------------
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.ListFormat.RemoveNumbers()
$p.Style = "Normal"
$p.Format.Alignment = 1
$p.Range.Text = "04/04 " + $endash + " 08/04"

$rDate = $d.Paragraphs.Last.Range
$rDate.Collapse(1)
$rDate.InsertParagraphBefore()

$n = $d.Paragraphs.Count
$headingIdx = $n - 1
$heading = $d.Paragraphs.Item($headingIdx)
$heading.Range.Font.Bold = 1
$heading.Range.Font.Size = 16
$heading.Range.Text = "Semaine 15"

# ---------------------------------------------------------------
# 6) "Reprise sur le front, fil d'arianne. Simplifie, et plus stable." (ilvl = 0)
# ---------------------------------------------------------------
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = "Paragraphedeliste"
$p.Range.ListFormat.ListLevelNumber = 1
$p.Range.Text = "Reprise sur le front, fil d" + $rsquo + "arianne. Simplifi" + [char]0x00E9 + ", et plus stable."

# ---------------------------------------------------------------
# 7) "Creation de la vue du details d'une instance. " (ilvl = 0, trailing space)
# ---------------------------------------------------------------
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.ListFormat.ListLevelNumber = 1
$p.Range.Text = "Creation de la vue du details d" + $rsquo + "une instance. "

# ---------------------------------------------------------------
# 8) "Maquettage, qui a ete ensuite valide ..." (ilvl = 1)
# ---------------------------------------------------------------
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.ListFormat.ListLevelNumber = 2
$p.Range.Text = "Maquettage, qui a " + [char]0x00E9 + "t" + [char]0x00E9 + " ensuite valid" + [char]0x00E9 + " par le restant de l" + $rsquo + "equipe, et modifi" + [char]0x00E9 + " suite au debat."

# ---------------------------------------------------------------
# 9) "Developpement de cette vue." (ilvl = 1)
# ---------------------------------------------------------------
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.ListFormat.ListLevelNumber = 2
$p.Range.Text = "Developpement de cette vue."

# ---------------------------------------------------------------
# 10) trailing empty paragraph
# ---------------------------------------------------------------
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.ListFormat.RemoveNumbers()
$p.Style = "Normal"

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
